$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: angle values (radians), i * pi/10 for i = 1..20, repeated twice
# across columns B..U and V..AO (40 columns total).
$angles = @(
    0.31415926535897898,
    0.62831853071795896,
    0.94247779607693805,
    1.2566370614359199,
    1.5707963267949001,
    1.8849555921538801,
    2.1991148575128601,
    2.5132741228718301,
    2.8274333882308098,
    3.14159265358979,
    3.4557519189487702,
    3.76991118430775,
    4.0840704496667302,
    4.3982297150257104,
    4.7123889803846897,
    5.0265482457436699,
    5.3407075111026501,
    5.6548667764616303,
    5.9690260418206096,
    6.2831853071795898
)

# Row 3: 0/1 flags, 40 values across columns B..AO
$flags = @(
    1,1,1,1,0,1,1,1,0,0,1,1,1,0,0,1,0,0,0,0,
    1,1,1,1,0,1,1,1,0,0,1,1,1,0,0,1,0,0,0,0
)

# Columns B (2) .. AO (41) -> 40 columns
for ($i = 0; $i -lt 40; $i++) {
    $col = $i + 2
    $ws.Cells.Item(2, $col).Value = $angles[$i % 20]
    $ws.Cells.Item(3, $col).Value = $flags[$i]
}

[void]$ws.Range("J17").Select()
